$d = $word.ActiveDocument
$wns = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently sits right after "Minutes week 8".
#    It needs to move to the end of the (new) "Chairman" line, so drop it
#    here and re-add it later once that line exists.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2. Insert three brand-new meta-data lines right after the title, re-using
#    the simple (un-styled) paragraph look of the existing meta-data lines.
# ---------------------------------------------------------------------------
$anchor = $d.Paragraphs.Item(2).Range
$anchor.InsertParagraphBefore()
$anchor.InsertParagraphBefore()
$anchor.InsertParagraphBefore()

$dateXml = '<w:p' + $wns + '><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>Date and time:</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> 23-04-2018, 13:45</w:t></w:r>' + `
  '</w:p>'
$d.Paragraphs.Item(2).Range.InsertXML($dateXml) | Out-Null

$locationXml = '<w:p' + $wns + '><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>Location</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>: Open space, Fontys R1</w:t></w:r>' + `
  '</w:p>'
$d.Paragraphs.Item(3).Range.InsertXML($locationXml) | Out-Null

$chairmanXml = '<w:p' + $wns + '><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>Chairman</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Rostislav</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Tinchev</w:t></w:r>' + `
  '</w:p>'
$d.Paragraphs.Item(4).Range.InsertXML($chairmanXml) | Out-Null

# Re-anchor "_GoBack" at the end of the new Chairman line (zero length,
# right before the paragraph mark -- same placement it had originally).
$chairmanRange = $d.Paragraphs.Item(4).Range
$goBackPoint = $d.Range($chairmanRange.End - 1, $chairmanRange.End - 1)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null

# ---------------------------------------------------------------------------
# 3. Drop the old, fragmented "Date and time" / "Location" / "Chairman"
#    paragraphs -- they have now been superseded by the ones just added.
#    After the inserts above they live at indexes 5, 6 and 7.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(5).Range.Delete() | Out-Null
$d.Paragraphs.Item(5).Range.Delete() | Out-Null
$d.Paragraphs.Item(5).Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 4. Expand the "Attendees" line with the full roster. Keep the existing
#    bold "Attendees" label run untouched and only replace what follows it.
# ---------------------------------------------------------------------------
$attendees = $d.Paragraphs.Item(6)
$labelEnd = $attendees.Range.Start + 9   # length of "Attendees"
$rest = $d.Range($labelEnd, $attendees.Range.End - 1)

$attendeesXml = '<w:p' + $wns + '>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Chung </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t>Kuah</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, Monika Kerulyte, Ignas Kybransas, </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:bCs/><w:color w:val="1D2129"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Teodor Genov</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:bCs/><w:color w:val="1D2129"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Yoanna Borisova</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:bCs/><w:color w:val="1D2129"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Rostislav Tinchev</w:t></w:r>' + `
  '</w:p>'
$rest.InsertXML($attendeesXml) | Out-Null

# ---------------------------------------------------------------------------
# 5. Merge the two "Meeting duration" trailer runs ("  " + "45mins") into
#    a single run, keeping the paragraph's own identity untouched.
# ---------------------------------------------------------------------------
$durationPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$durFull = $durationPara.Range
$labelLen = 17   # length of "Meeting duration:"
$fortyFive = $d.Range($durFull.Start + $labelLen + 2, $durFull.End - 1)
$fortyFive.InsertBefore("  ")
$spaceRun = $d.Range($durFull.Start + $labelLen, $durFull.Start + $labelLen + 2)
$spaceRun.Delete() | Out-Null

Write-Host "done"
